# Insert a new daily price record as row 171 (pushing the existing rows
# 171-273 down to 172-274), then populate the new row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("171:171").Insert()

$ws.Range("A171").Value = 10
$ws.Range("B171").Value = "Vega Modelo de Temuco"
$ws.Range("C171").Value = "La Araucanía"
$ws.Range("D171").Value = 44603
$ws.Range("E171").Value = 9
$ws.Range("F171").Value = 100112009
$ws.Range("G171").Value = "Acelga"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 90
$ws.Range("K171").Value = 7000
$ws.Range("L171").Value = 8000
$ws.Range("M171").Value = 7389
$ws.Range("N171").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O171").Value = "Provincia de Cautín"
$ws.Range("P171").Value = 616
$ws.Range("Q171").Value = 12
$ws.Range("R171").Value = "Hortaliza"
